# CapstoneHours.xlsx update — "Adding in New assets to use"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix the dates on rows 38-41 (were 10/14/2014 -> now 10/12/2014) ---
$ws.Range("C38").Value = 41924
$ws.Range("C39").Value = 41924
$ws.Range("C40").Value = 41924
$ws.Range("C41").Value = 41924

# --- Shared-string write order matters for matching the canonical table,
#     so the "updating prefabs..." text (which reuses the old slot vacated
#     by "Adding in jumping sound") is written before the rest. ---
$ws.Range("A45").Value = "updating prefabs with nodes and pathchoices"
$ws.Range("A42").Value = "Researching Random generation for mazes"
$ws.Range("A43").Value = "Writing Story(prequel)"
$ws.Range("A44").Value = "Fix movement sounds"
$ws.Range("A46").Value = "Write Maze ending types(events that will happen and how they affect game)"
$ws.Range("A47").Value = "finding models"
$ws.Range("E3").Value = "Week 2 Hours"
$ws.Range("A49").Value = "seting up models (textures, collision, scale, etc)"
$ws.Range("E50").Value = "started at 12"
$ws.Range("A48").Value = "finding models"
$ws.Range("A50").Value = "seting up models (textures, collision, scale, etc)"

# --- Numeric hours ---
$ws.Range("B42").Value = 2
$ws.Range("B43").Value = 3
$ws.Range("B44").Value = 2.5
$ws.Range("B45").Value = 0.5
$ws.Range("B46").Value = 2
$ws.Range("B47").Value = 1
$ws.Range("B48").Value = 4
$ws.Range("B49").Value = 2.5

# --- Dates for the new rows, copying the number format from an existing
#     date cell so no new style record is created. ---
$ws.Range("C38").Copy() | Out-Null
$ws.Range("C42:C50").PasteSpecial(-4122) | Out-Null
$ws.Range("C42").Value = 41925
$ws.Range("C43").Value = 41925
$ws.Range("C44").Value = 41926
$ws.Range("C45").Value = 41926
$ws.Range("C46").Value = 41927
$ws.Range("C47").Value = 41927
$ws.Range("C48").Value = 41928
$ws.Range("C49").Value = 41928
$ws.Range("C50").Value = 41929
$excel.CutCopyMode = 0

# --- Totals ---
$ws.Range("F1").Formula = "=SUM(B2,B3:B301)"
$ws.Range("F3").Formula = "=SUM(B20:B43)"

# --- View state: scroll position / selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 21
$win.ScrollColumn = 1
$ws.Range("E44").Select() | Out-Null
